# Insert a new data row at row 78 (pushing existing rows 78:161 down to 79:162)
# and populate the new row with the new price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 78. Excel will shift all rows
# 78..161 down to 79..162 and copy formatting (incl. the date number format
# used in column D) from the row above, just like a normal Excel insert.
$ws.Rows.Item(78).Insert()

# Populate the newly inserted row 78 with the new record.
$ws.Range("A78").Value = 5
$ws.Range("B78").Value = "Macroferia Regional de Talca"
$ws.Range("C78").Value = "Maule"
$ws.Range("D78").Value = 44494
$ws.Range("E78").Value = 7
$ws.Range("F78").Value = "Fruta"
$ws.Range("G78").Value = 100108
$ws.Range("H78").Value = "Tropicales y subtropicales"
$ws.Range("I78").Value = 100108005
$ws.Range("J78").Value = "Piña"
$ws.Range("K78").Value = "Caramelo"
$ws.Range("L78").Value = "Tercera"
$ws.Range("M78").Value = 210
$ws.Range("N78").Value = 19000
$ws.Range("O78").Value = 19000
$ws.Range("P78").Value = 19000
$ws.Range("Q78").Value = "`$/caja 16 unidades"
$ws.Range("R78").Value = "Ecuador"
$ws.Range("S78").Value = 1188
$ws.Range("T78").Value = 16
